# Regenerate save_data to use K (strike count / option count proxy) instead of
# the old "Strike#" derived value, and recompute std/mean-derived s_vals for
# the K column (column G) across all data rows.
#
# Target K values per row, keyed by worksheet row number, as recomputed by
# the upstream regeneration pipeline.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 1
    3 = 4
    4 = 2
    5 = 0
    6 = 2
    8 = 1
    9 = 1
    10 = 1
    11 = 2
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 2
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 2
    23 = 0
    24 = 1
    25 = 1
    27 = 1
    28 = 2
    29 = 2
    30 = 2
    31 = 0
    32 = 0
    33 = 1
    34 = 1
    35 = 1
    36 = 1
    37 = 1
    38 = 0
    40 = 3
    41 = 1
    42 = 1
    43 = 2
    44 = 1
    45 = 2
    46 = 2
    47 = 2
    48 = 2
    49 = 1
    51 = 2
    52 = 0
    53 = 2
    54 = 0
    55 = 3
    56 = 1
    57 = 1
    58 = 3
    59 = 2
    60 = 6
    61 = 2
    62 = 2
    63 = 2
    64 = 5
    65 = 2
    66 = 0
    67 = 1
    68 = 0
    69 = 4
    70 = 3
    71 = 3
    72 = 0
    73 = 3
    74 = 1
    75 = 1
    76 = 1
    77 = 1
    78 = 0
    79 = 0
    80 = 2
    81 = 2
    82 = 2
    83 = 2
    84 = 1
    87 = 2
    88 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
